$wb = $excel.ActiveWorkbook

# "Grafo 1" keeps its same structure, but the person formerly called
# "Elmago" is now called "Mariano" (a shared-string rename).
$ws1 = $wb.Worksheets.Item("Grafo 1")
[void]$ws1.Cells.Replace("Elmago", "Mariano")

# "Grafo 2" gets a brand-new set of rows describing a different little
# graph of people (Juan, Pedro, Ana, Luis, Hugo) and their relations.
$ws2 = $wb.Worksheets.Item("Grafo 2")

# Every data cell in this sheet uses the same cell style (index 1, the
# sheet's standard Arial style). Copy that style onto each target cell
# before writing its value so newly-created cells match the existing ones.
$targetCells = @("A1","B1","C1","D1","E1","A2","B2","C2","A3","B3","C3","D3","E3","A4","B4","C4","A5")
foreach ($addr in $targetCells) {
    $ws2.Range("A1").Copy()
    $ws2.Range($addr).PasteSpecial(-4122)
}

# Row 1: Juan -> Pedro (compañero), Juan -> Ana (amigo personal)
$ws2.Range("A1").Value = "Juan"
$ws2.Range("B1").Value = "Pedro"
$ws2.Range("C1").Value = "compañero"
$ws2.Range("D1").Value = "Ana"
$ws2.Range("E1").Value = "amigo personal"

# Row 2: Pedro -> Luis (amigo personal)
$ws2.Range("A2").Value = "Pedro"
$ws2.Range("B2").Value = "Luis"
$ws2.Range("C2").Value = "amigo personal"

# Row 3: Ana -> Hugo (compañero), Ana -> Luis (conocido)
$ws2.Range("A3").Value = "Ana"
$ws2.Range("B3").Value = "Hugo"
$ws2.Range("C3").Value = "compañero"
$ws2.Range("D3").Value = "Luis"
$ws2.Range("E3").Value = "conocido"

# Row 4: Luis -> Hugo (amigo personal)
$ws2.Range("A4").Value = "Luis"
$ws2.Range("B4").Value = "Hugo"
$ws2.Range("C4").Value = "amigo personal"

# Row 5: Hugo (no outgoing edge listed)
$ws2.Range("A5").Value = "Hugo"
